$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.884.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = "'3.498.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'594.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = "'0.132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.88%  '
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").Value = "'0.432"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = "'4.100.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = "'29.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.76%  '
$ws.Range("D15").Value = "'66.880.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = "'3.500.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = "'14.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '
$ws.Range("D20").Value = "'394.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").Value = "'7.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = "'73.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = "'10.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").Value = "'7.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("D35").Value = "'162.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("D36").Value = "'0.879"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").Value = "'4.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = "'0.0738"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").Value = "'2.834.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.39%  '
$ws.Range("D42").Value = "'27.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").Value = "'26.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("D44").Value = "'42.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("E45").Value = '  +2.18%  '
$ws.Range("E46").Value = '  -3.01%  '
$ws.Range("D47").Value = "'337.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("D48").Value = "'34.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").Value = "'0.842"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.04%  '
